$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metro_budget")

# ---------------------------------------------------------------------------
# Helper: set a cell's value/formula then copy number-format/style from a
# template cell that already carries the desired style index, without
# disturbing the template cell itself.
# ---------------------------------------------------------------------------
function Set-StyledValue($rng, $value, $styleSrc) {
    $rng.Value = $value
    $styleSrc.Copy()
    $rng.PasteSpecial(-4122)
}

function Set-StyledFormula($rng, $formula, $styleSrc) {
    $rng.Formula = $formula
    if ($styleSrc -ne $null) {
        $styleSrc.Copy()
        $rng.PasteSpecial(-4122)
    }
}

# Style template cells (already present in the sheet, untouched by copy):
$styleA55  = $ws.Range("A55")   # s="1"
$styleE89  = $ws.Range("E89")   # s="10"
$styleD55  = $ws.Range("D55")   # s="14"
$styleE1   = $ws.Range("E1")    # s="9"

# ---------------------------------------------------------------------------
# Row 54: new "Question 7" header cell in F54
# ---------------------------------------------------------------------------
Set-StyledValue $ws.Range("F54") "Question 7" $styleE89

# ---------------------------------------------------------------------------
# Row 55: column headers for the new Question 7 block (F:H)
# ---------------------------------------------------------------------------
Set-StyledValue $ws.Range("F55") "FY17_diff" $styleA55
Set-StyledValue $ws.Range("G55") "FY18_diff" $styleA55
Set-StyledValue $ws.Range("H55") "FY19_diff" $styleD55

# ---------------------------------------------------------------------------
# Rows 56-61: VLOOKUP exercise mirrored into F:H
# ---------------------------------------------------------------------------
Set-StyledFormula $ws.Range("H56") "=INDEX(`$N`$2:`$N`$52,MATCH(A56,`$A`$2:`$A`$52,0))" $styleE1

Set-StyledFormula $ws.Range("F57") "=INDEX(`$D`$2:`$D`$52, MATCH(A58,`$A`$2:`$A`$52,0))" $null
Set-StyledFormula $ws.Range("G57") "=INDEX(`$I`$2:`$I`$52,MATCH(A57,`$A`$2:`$A`$52,0))" $styleE1
Set-StyledFormula $ws.Range("H57") "=INDEX(`$N`$2:`$N`$52,MATCH(A57,`$A`$2:`$A`$52,0))" $styleE1

Set-StyledFormula $ws.Range("G58") "=INDEX(`$I`$2:`$I`$52,MATCH(A58,`$A`$2:`$A`$52,0))" $styleE1
Set-StyledFormula $ws.Range("H58") "=INDEX(`$N`$2:`$N`$52,MATCH(A58,`$A`$2:`$A`$52,0))" $styleE1

Set-StyledFormula $ws.Range("F59") "=INDEX(`$D`$2:`$D`$52, MATCH(A59,`$A`$2:`$A`$52,0))" $null
Set-StyledFormula $ws.Range("G59") "=INDEX(`$I`$2:`$I`$52,MATCH(A59,`$A`$2:`$A`$52,0))" $styleE1
Set-StyledFormula $ws.Range("H59") "=INDEX(`$N`$2:`$N`$52,MATCH(A59,`$A`$2:`$A`$52,0))" $styleE1

Set-StyledFormula $ws.Range("F60") "=INDEX(`$D`$2:`$D`$52, MATCH(A60,`$A`$2:`$A`$52,0))" $null
Set-StyledFormula $ws.Range("G60") "=INDEX(`$I`$2:`$I`$52,MATCH(A60,`$A`$2:`$A`$52,0))" $styleE1
Set-StyledFormula $ws.Range("H60") "=INDEX(`$N`$2:`$N`$52,MATCH(A60,`$A`$2:`$A`$52,0))" $styleE1

Set-StyledFormula $ws.Range("F61") "=INDEX(`$D`$2:`$D`$52, MATCH(A61,`$A`$2:`$A`$52,0))" $null
Set-StyledFormula $ws.Range("G61") "=INDEX(`$I`$2:`$I`$52,MATCH(A61,`$A`$2:`$A`$52,0))" $styleE1
Set-StyledFormula $ws.Range("H61") "=INDEX(`$N`$2:`$N`$52,MATCH(A61,`$A`$2:`$A`$52,0))" $styleE1

# Row 62: stray blank cell F62 (present in the target file, no content)
$ws.Range("F62").Value = ""

# ---------------------------------------------------------------------------
# Row 64: column headers for Question 8 (XLOOKUP) block, add F:H headers
# ---------------------------------------------------------------------------
Set-StyledValue $ws.Range("F64") "FY17_diff" $styleA55
Set-StyledValue $ws.Range("G64") "FY18_diff" $styleA55
Set-StyledValue $ws.Range("H64") "FY19_diff" $styleD55

# ---------------------------------------------------------------------------
# Row 73: column headers for Question 9 (INDEX/MATCH) block, add F:H headers
# ---------------------------------------------------------------------------
Set-StyledValue $ws.Range("F73") "FY17_diff" $styleA55
Set-StyledValue $ws.Range("G73") "FY18_diff" $styleA55
Set-StyledValue $ws.Range("H73") "FY19_diff" $styleD55

# ---------------------------------------------------------------------------
# Rows 74-79: tweak existing INDEX/MATCH formulas ($A7x anchors, C column
# switched from H to I)
# ---------------------------------------------------------------------------
$ws.Range("B74").Formula = "=INDEX(`$D`$2:`$D`$52, MATCH(`$A`74,`$A`$2:`$A`$52,0))"
$ws.Range("C74").Formula = "=INDEX(`$I`$2:`$I`$52,MATCH(`$A`74,`$A`$2:`$A`$52,0))"
$ws.Range("D74").Formula = "=INDEX(`$N`$2:`$N`$52,MATCH(`$A`74,`$A`$2:`$A`$52,0))"

$ws.Range("B75").Formula = "=INDEX(`$D`$2:`$D`$52, MATCH(`$A`75,`$A`$2:`$A`$52,0))"
$ws.Range("C75").Formula = "=INDEX(`$I`$2:`$I`$52,MATCH(`$A`75,`$A`$2:`$A`$52,0))"
$ws.Range("D75").Formula = "=INDEX(`$N`$2:`$N`$52,MATCH(`$A`75,`$A`$2:`$A`$52,0))"

$ws.Range("B76").Formula = "=INDEX(`$D`$2:`$D`$52, MATCH(`$A`76,`$A`$2:`$A`$52,0))"
$ws.Range("C76").Formula = "=INDEX(`$I`$2:`$I`$52,MATCH(`$A`76,`$A`$2:`$A`$52,0))"
$ws.Range("D76").Formula = "=INDEX(`$N`$2:`$N`$52,MATCH(`$A`76,`$A`$2:`$A`$52,0))"

$ws.Range("B77").Formula = "=INDEX(`$D`$2:`$D`$52, MATCH(`$A`77,`$A`$2:`$A`$52,0))"
$ws.Range("C77").Formula = "=INDEX(`$I`$2:`$I`$52,MATCH(`$A`77,`$A`$2:`$A`$52,0))"
$ws.Range("D77").Formula = "=INDEX(`$N`$2:`$N`$52,MATCH(`$A`77,`$A`$2:`$A`$52,0))"

$ws.Range("B78").Formula = "=INDEX(`$D`$2:`$D`$52, MATCH(`$A`78,`$A`$2:`$A`$52,0))"
$ws.Range("C78").Formula = "=INDEX(`$I`$2:`$I`$52,MATCH(`$A`78,`$A`$2:`$A`$52,0))"
$ws.Range("D78").Formula = "=INDEX(`$N`$2:`$N`$52,MATCH(`$A`78,`$A`$2:`$A`$52,0))"

$ws.Range("B79").Formula = "=INDEX(`$D`$2:`$D`$52, MATCH(`$A`79,`$A`$2:`$A`$52,0))"
$ws.Range("C79").Formula = "=INDEX(`$I`$2:`$I`$52,MATCH(`$A`79,`$A`$2:`$A`$52,0))"
$ws.Range("D79").Formula = "=INDEX(`$N`$2:`$N`$52,MATCH(`$A`79,`$A`$2:`$A`$52,0))"

# ---------------------------------------------------------------------------
# Row 87: pick a different department for the summary chart (was "Community
# Education Commission", now "Arts Commission")
# ---------------------------------------------------------------------------
$ws.Range("B87").Value = "Arts Commission"

Write-Host "sheet1 data edits complete"
